$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header-row cells (row 1), added above the title row.
# "32423" looks numeric, so force it to be stored as text (leading
# apostrophe), then reset the cell style back to Normal so it doesn't
# pick up a lingering quote-prefix / number-format style.
$ws.Range("C1").Value = "dadaw"
$ws.Range("E1").Value = "'32423"
$ws.Range("E1").Style = "Normal"

# Fall 2022 / Spring 2022 / Summer 2022 block (rows 4-10)
# A new course (PSYC 1101) is inserted at the top, shifting the existing
# Fall 2022 courses down by one row; two more rows are appended at the
# bottom (CPSC 3415, CPSC 4000) and a new Summer 2022 course (CPSC 4899)
# is added next to the first row.
$ws.Range("A4").Value = "PSYC 1101"
$ws.Range("B4").Value = 3
$ws.Range("E4").Value = "CPSC 4899"
$ws.Range("F4").Value = 3

$ws.Range("A5").Value = "POLS 1101"
$ws.Range("B5").Value = 3

$ws.Range("A6").Value = "PSYC 1105"
$ws.Range("B6").Value = 2

$ws.Range("A7").Value = "DSCI 3111"
$ws.Range("B7").Value = 3

$ws.Range("A8").Value = "CPSC 3121"
$ws.Range("B8").Value = 3

$ws.Range("A9").Value = "CPSC 3415"
$ws.Range("B9").Value = 1

$ws.Range("A10").Value = "CPSC 4000"
$ws.Range("B10").Value = 0

# Fall 2023 block (rows 13-19) loses its last two entries
# (CPSC 4205, CPSC 4555) leaving just CPSC 4175 / CPSC 4176.
$ws.Range("A14").ClearContents()
$ws.Range("B14").ClearContents()
$ws.Range("A15").ClearContents()
$ws.Range("B15").ClearContents()

$wb.Save()
